$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 193, shifting the existing rows 193-196 down to 194-197.
$ws.Rows.Item(193).Insert()

# Populate the newly inserted row 193 with the new weekly record.
$ws.Range("A193").Value = 9
$ws.Range("B193").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C193").Value = "Metropolitana"
$ws.Range("D193").Value = 45239
$ws.Range("E193").Value = 13
$ws.Range("F193").Value = 100112022
$ws.Range("G193").Value = "Arveja Verde"
$ws.Range("H193").Value = "Sin especificar"
$ws.Range("I193").Value = "Primera"
$ws.Range("J193").Value = 52
$ws.Range("K193").Value = 18000
$ws.Range("L193").Value = 19000
$ws.Range("M193").Value = 18500
$ws.Range("N193").Value = "$/saco 25 kilos"
$ws.Range("O193").Value = "Provincia de Huasco"
$ws.Range("P193").Value = 740
$ws.Range("Q193").Value = 25
$ws.Range("R193").Value = "Hortaliza"
